$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(132, 8).Value = 1768.0303
$ws.Cells.Item(132, 9).Value = 1355.4755
$ws.Cells.Item(132, 11).Value = 4066.4265
$ws.Cells.Item(132, 13).Value = -1536.4265
$ws.Cells.Item(134, 8).Value = 56780
$ws.Cells.Item(134, 10).Value = 56780
$ws.Cells.Item(134, 12).Value = 56780
$ws.Cells.Item(134, 14).Value = -66920
$ws.Cells.Item(137, 8).Value = 1235.0264
$ws.Cells.Item(137, 9).Value = 1253.5385
$ws.Cells.Item(137, 10).Value = 1194.9166
$ws.Cells.Item(137, 11).Value = 3760.6155
$ws.Cells.Item(137, 12).Value = 3584.7498
$ws.Cells.Item(137, 13).Value = -1210.6155
$ws.Cells.Item(137, 14).Value = -8684.7498
$ws.Cells.Item(138, 8).Value = 4618.7
$ws.Cells.Item(138, 9).Value = 1986.85
$ws.Cells.Item(138, 10).Value = 5495.9834
$ws.Cells.Item(138, 11).Value = 5960.549999999999
$ws.Cells.Item(138, 12).Value = 16487.9502
$ws.Cells.Item(138, 13).Value = -820.5499999999993
$ws.Cells.Item(138, 14).Value = -26767.9502
$ws.Cells.Item(140, 8).Value = 69880
$ws.Cells.Item(140, 10).Value = 69880
$ws.Cells.Item(140, 12).Value = 69880
$ws.Cells.Item(140, 14).Value = -80240
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(4, 8).Value = 250
$ws.Cells.Item(4, 9).Value = 250
$ws.Cells.Item(4, 11).Value = 250
$ws.Cells.Item(4, 13).Value = -134
$ws.Cells.Item(5, 8).Value = 475.25
$ws.Cells.Item(5, 9).Value = 467
$ws.Cells.Item(5, 10).Value = 500
$ws.Cells.Item(5, 11).Value = 467
$ws.Cells.Item(5, 12).Value = 500
$ws.Cells.Item(5, 13).Value = -355
$ws.Cells.Item(5, 14).Value = -724
$ws.Cells.Item(61, 8).Value = 248978.83
$ws.Cells.Item(61, 9).Value = 6478.6523
$ws.Cells.Item(61, 10).Value = 558840.2
$ws.Cells.Item(61, 11).Value = 6478.6523
$ws.Cells.Item(61, 12).Value = 558840.2
$ws.Cells.Item(61, 13).Value = -6266.6523
$ws.Cells.Item(61, 14).Value = -559264.2
$ws.Cells.Item(102, 8).Value = 2850290.8
$ws.Cells.Item(102, 9).Value = 2850290.8
$ws.Cells.Item(102, 11).Value = 2850290.8
$ws.Cells.Item(102, 13).Value = -2848668.8
$ws.Cells.Item(136, 8).Value = 248978.83
$ws.Cells.Item(136, 9).Value = 6478.6523
$ws.Cells.Item(136, 10).Value = 558840.2
$ws.Cells.Item(136, 11).Value = 19435.9569
$ws.Cells.Item(136, 12).Value = 1676520.6
$ws.Cells.Item(136, 13).Value = -16885.9569
$ws.Cells.Item(136, 14).Value = -1681620.6
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 475.25
$ws.Cells.Item(4, 9).Value = 467
$ws.Cells.Item(4, 10).Value = 500
$ws.Cells.Item(4, 11).Value = 467
$ws.Cells.Item(4, 12).Value = 500
$ws.Cells.Item(4, 13).Value = -352
$ws.Cells.Item(4, 14).Value = -730
$ws.Cells.Item(22, 8).Value = 593.3
$ws.Cells.Item(22, 9).Value = 658.5
$ws.Cells.Item(22, 10).Value = 549.8333
$ws.Cells.Item(22, 11).Value = 658.5
$ws.Cells.Item(22, 12).Value = 549.8333
$ws.Cells.Item(22, 13).Value = -485.5
$ws.Cells.Item(22, 14).Value = -895.8333
$ws.Cells.Item(26, 8).Value = 14000
$ws.Cells.Item(26, 9).Value = 14000
$ws.Cells.Item(26, 10).Value = 0
$ws.Cells.Item(26, 11).Value = 14000
$ws.Cells.Item(26, 12).Value = 0
$ws.Cells.Item(26, 14).Value = -13708
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(122, 8).Value = 4023.65
$ws.Cells.Item(122, 9).Value = 3058.7856
$ws.Cells.Item(122, 10).Value = 6275
$ws.Cells.Item(122, 11).Value = 9176.356800000001
$ws.Cells.Item(122, 12).Value = 18825
$ws.Cells.Item(122, 13).Value = -6726.356800000001
$ws.Cells.Item(122, 14).Value = -23725
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(38, 8).Value = 4166944.8
$ws.Cells.Item(38, 9).Value = 7143227
$ws.Cells.Item(38, 10).Value = 149.3
$ws.Cells.Item(38, 11).Value = 21429681
$ws.Cells.Item(38, 12).Value = 447.9
$ws.Cells.Item(38, 13).Value = -21429334
$ws.Cells.Item(38, 14).Value = -1141.9
$ws.Cells.Item(44, 8).Value = 804.95
$ws.Cells.Item(44, 10).Value = 1091.5834
$ws.Cells.Item(44, 12).Value = 3274.7502
$ws.Cells.Item(44, 14).Value = -4070.7502
$ws.Cells.Item(113, 8).Value = 222717.72
$ws.Cells.Item(113, 9).Value = 505.25
$ws.Cells.Item(113, 10).Value = 303522.25
$ws.Cells.Item(113, 11).Value = 1515.75
$ws.Cells.Item(113, 12).Value = 910566.75
$ws.Cells.Item(113, 13).Value = 654.25
$ws.Cells.Item(113, 14).Value = -914906.75
$ws.Cells.Item(132, 8).Value = 3576.8125
$ws.Cells.Item(132, 9).Value = 2526
$ws.Cells.Item(132, 10).Value = 3927.0833
$ws.Cells.Item(132, 11).Value = 22734
$ws.Cells.Item(132, 12).Value = 35343.7497
$ws.Cells.Item(132, 13).Value = -20204
$ws.Cells.Item(132, 14).Value = -40403.7497
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(11, 8).Value = 12092455
$ws.Cells.Item(11, 9).Value = 11000625
$ws.Cells.Item(11, 10).Value = 15004000
$ws.Cells.Item(11, 11).Value = 11000625
$ws.Cells.Item(11, 12).Value = 15004000
$ws.Cells.Item(11, 13).Value = -11000486
$ws.Cells.Item(11, 14).Value = -15004278
$ws.Cells.Item(20, 8).Value = 4006653.5
$ws.Cells.Item(20, 9).Value = 10004002
$ws.Cells.Item(20, 11).Value = 10004002
$ws.Cells.Item(20, 13).Value = -10003757
$ws.Cells.Item(122, 8).Value = 119531870
$ws.Cells.Item(122, 10).Value = 41671636
$ws.Cells.Item(122, 12).Value = 125014908
$ws.Cells.Item(122, 14).Value = -125019808
$ws.Cells.Item(136, 8).Value = 52833.332
$ws.Cells.Item(136, 10).Value = 52833.332
$ws.Cells.Item(136, 12).Value = 158499.996
$ws.Cells.Item(136, 14).Value = -163599.996
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1032.1111
$ws.Cells.Item(22, 9).Value = 1137.8
$ws.Cells.Item(22, 10).Value = 900
$ws.Cells.Item(22, 11).Value = 1137.8
$ws.Cells.Item(22, 12).Value = 900
$ws.Cells.Item(22, 13).Value = -842.8
$ws.Cells.Item(22, 14).Value = -1490
$ws.Cells.Item(24, 8).Value = 50007
$ws.Cells.Item(24, 10).Value = 50007
$ws.Cells.Item(24, 12).Value = 50007
$ws.Cells.Item(24, 14).Value = -50693
$ws.Cells.Item(25, 8).Value = 8625
$ws.Cells.Item(25, 9).Value = 1500
$ws.Cells.Item(25, 10).Value = 30000
$ws.Cells.Item(25, 11).Value = 1500
$ws.Cells.Item(25, 12).Value = 30000
$ws.Cells.Item(25, 13).Value = -1270
$ws.Cells.Item(25, 14).Value = -30460
$ws.Cells.Item(27, 8).Value = 1032.1111
$ws.Cells.Item(27, 9).Value = 1137.8
$ws.Cells.Item(27, 10).Value = 900
$ws.Cells.Item(27, 11).Value = 1137.8
$ws.Cells.Item(27, 12).Value = 900
$ws.Cells.Item(27, 13).Value = -1030.8
$ws.Cells.Item(27, 14).Value = -1114
$ws.Cells.Item(40, 8).Value = 2110.9583
$ws.Cells.Item(40, 9).Value = 1969.8889
$ws.Cells.Item(40, 10).Value = 2534.1667
$ws.Cells.Item(40, 11).Value = 1969.8889
$ws.Cells.Item(40, 12).Value = 2534.1667
$ws.Cells.Item(40, 13).Value = -1833.8889
$ws.Cells.Item(40, 14).Value = -2806.1667
$ws.Cells.Item(46, 8).Value = 1231.8182
$ws.Cells.Item(46, 9).Value = 962.5
$ws.Cells.Item(46, 10).Value = 1385.7142
$ws.Cells.Item(46, 11).Value = 962.5
$ws.Cells.Item(46, 12).Value = 1385.7142
$ws.Cells.Item(46, 13).Value = -774.5
$ws.Cells.Item(46, 14).Value = -1761.7142
$ws.Cells.Item(74, 8).Value = 21574
$ws.Cells.Item(74, 10).Value = 21574
$ws.Cells.Item(74, 12).Value = 21574
$ws.Cells.Item(74, 14).Value = -23570
$ws.Cells.Item(77, 8).Value = 21574
$ws.Cells.Item(77, 10).Value = 21574
$ws.Cells.Item(77, 12).Value = 64722
$ws.Cells.Item(77, 14).Value = -74706
$ws.Cells.Item(132, 8).Value = 5856.4546
$ws.Cells.Item(132, 9).Value = 6257.794
$ws.Cells.Item(132, 10).Value = 4491.9
$ws.Cells.Item(132, 11).Value = 18773.382
$ws.Cells.Item(132, 12).Value = 13475.7
$ws.Cells.Item(132, 13).Value = -16243.382
$ws.Cells.Item(132, 14).Value = -18535.7
$ws.Cells.Item(136, 8).Value = 10322.4375
$ws.Cells.Item(136, 9).Value = 6674.75
$ws.Cells.Item(136, 10).Value = 21265.5
$ws.Cells.Item(136, 11).Value = 20024.25
$ws.Cells.Item(136, 12).Value = 63796.5
$ws.Cells.Item(136, 13).Value = -17474.25
$ws.Cells.Item(136, 14).Value = -68896.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(20, 8).Value = 10000000
$ws.Cells.Item(20, 9).Value = 10000000
$ws.Cells.Item(20, 11).Value = 10000000
$ws.Cells.Item(20, 13).Value = -9999760
$ws.Cells.Item(132, 8).Value = 1178.36
$ws.Cells.Item(132, 9).Value = 558.9706
$ws.Cells.Item(132, 11).Value = 1676.9118
$ws.Cells.Item(132, 13).Value = 853.0882000000001
$ws.Cells.Item(136, 8).Value = 2491.26
$ws.Cells.Item(136, 9).Value = 2458.6333
$ws.Cells.Item(136, 10).Value = 2540.2
$ws.Cells.Item(136, 11).Value = 7375.8999
$ws.Cells.Item(136, 12).Value = 7620.599999999999
$ws.Cells.Item(136, 13).Value = -4825.8999
$ws.Cells.Item(136, 14).Value = -12720.6

Write-Output "Applied 202 cell updates"